# Update workbook for "Add data for 2022-04-19" (dates advance from
# through-04-10 to through-04-11 in the sheet name/header, with one more
# day of data folded into the April row and the Total row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet (tab name + workbook.xml <sheet name=...>)
$ws.Name = "Through 2022-04-11"

# 2) Update the "April (through 04-10)" label cell (row 5, column A)
$ws.Range("A5").Value = "April (through 04-11)"

# 3) Update April row (row 5) counts for columns D..I (2017-2022)
$ws.Range("D5").Value = 21
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 31
$ws.Range("I5").Value = 40

# 4) Update Total row (row 6) counts for columns D..I (2017-2022)
$ws.Range("D6").Value = 210
$ws.Range("E6").Value = 216
$ws.Range("F6").Value = 128
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 454
$ws.Range("I6").Value = 474
